$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for each data row.
# Every data row (2 through 307) had this date bumped from 2023-09-03
# (serial 45172) to 2023-09-06 (serial 45175).
$newDate = Get-Date -Year 2023 -Month 9 -Day 6 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

for ($r = 2; $r -le 307; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
